# Remove redundant unit-related slots now that quantitative fields use
# QuantityValue (which already carries its own unit).
#
#   Sample.concentration_unit            -> delete column F on "Sample" sheet
#   StorageConditions.temperature_unit   -> delete column B on "StorageConditions" sheet
#
# Deleting the entire column shifts everything to its right one slot to the
# left, which also re-bases the sheet's used-range / dimension and drops any
# data validation whose sqref pointed only at the removed column.

$wb = $excel.ActiveWorkbook

$sampleSheet = $wb.Worksheets.Item("Sample")
$sampleSheet.Range("F1").EntireColumn.Delete()

$storageSheet = $wb.Worksheets.Item("StorageConditions")
$storageSheet.Range("B1").EntireColumn.Delete()
